# Daily attendance processing - 2026-01-09 16:40:12
# Normalizes the "Recorded By" (column G) cell values so that when a
# session was recorded by both the automated System and a named user,
# the two names are listed in a consistent order (swap the two
# comma-separated entries wherever "System" appears alongside another
# recorder).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$colRecordedBy = 7  # Column G: "Recorded By"
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colRecordedBy)
    $value = $cell.Value2

    if ($value -ne $null -and $value -like "*System*" -and $value -like "*,*") {
        $parts = $value -split ",\s*"
        if ($parts.Count -eq 2) {
            $swapped = $parts[1].Trim() + ", " + $parts[0].Trim()
            if ($swapped -ne $value) {
                $cell.Value = $swapped
                $changed++
            }
        }
    }
}

Write-Host "Swapped $changed cells"
